$d = $word.ActiveDocument

$replacements = @(
    @("50×77=3850", "89×25=2225"),
    @("60×93=5580", "13×12=156"),
    @("70×60=4200", "34×16=544"),
    @("69×48=3312", "89×73=6497"),
    @("95×66=6270", "68×22=1496"),
    @("33×17=561", "21×92=1932"),
    @("43×56=2408", "68×31=2108"),
    @("15×37=555", "22×20=440"),
    @("59×59=3481", "39×26=1014"),
    @("58×80=4640", "15×38=570"),
    @("30×42=1260", "80×49=3920"),
    @("72×34=2448", "20×86=1720"),
    @("14×80=1120", "80×48=3840"),
    @("64×14=896", "66×31=2046"),
    @("13×55=715", "43×72=3096"),
    @("90×77=6930", "81×46=3726"),
    @("17×91=1547", "87×14=1218"),
    @("98×40=3920", "47×91=4277"),
    @("34×20=680", "84×90=7560"),
    @("91×46=4186", "72×76=5472"),
    @("64×86=5504", "46×56=2576"),
    @("33×48=1584", "21×55=1155"),
    @("88×94=8272", "26×12=312"),
    @("74×86=6364", "25×54=1350"),
    @("81×86=6966", "24×73=1752")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
